$wb = $excel.ActiveWorkbook

# Hunk 0: sheet ALC, row 58
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1122.0769
$ws.Range("I58").Value = 200.875
$ws.Range("K58").Value = 602.625
$ws.Range("M58").Value = -452.625

# Hunk 1: sheet ALC, row 70
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 41917892
$ws.Range("J70").Value = 1795
$ws.Range("L70").Value = 5385
$ws.Range("N70").Value = -5925

# Hunk 2: sheet ALC, row 73
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 41917892
$ws.Range("J73").Value = 1795
$ws.Range("L73").Value = 5385
$ws.Range("N73").Value = -7257

# Hunk 3: sheet ALC, row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 2425.2856
$ws.Range("J86").Value = 2743.1667
$ws.Range("L86").Value = 2743.1667
$ws.Range("N86").Value = -4989.1667

# Hunk 4: sheet ALC, row 87
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 25390
$ws.Range("J87").Value = 25390
$ws.Range("L87").Value = 25390
$ws.Range("N87").Value = -27886

# Hunk 5: sheet ALC, row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 2425.2856
$ws.Range("J89").Value = 2743.1667
$ws.Range("L89").Value = 13715.8335
$ws.Range("N89").Value = -24947.8335

# Hunk 6: sheet ALC, row 90
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 25390
$ws.Range("J90").Value = 25390
$ws.Range("L90").Value = 76170
$ws.Range("N90").Value = -88650

# Hunk 7: sheet ALC, row 101
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 1668.0834
$ws.Range("I101").Value = 522.4
$ws.Range("J101").Value = 2486.4285
$ws.Range("K101").Value = 1567.2
$ws.Range("L101").Value = 7459.2855
$ws.Range("M101").Value = 54.80000000000018
$ws.Range("N101").Value = -10703.2855

# Hunk 8: sheet ALC, row 111
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 7755.45
$ws.Range("I111").Value = 2656
$ws.Range("J111").Value = 13988.111
$ws.Range("K111").Value = 7968
$ws.Range("L111").Value = 41964.333
$ws.Range("M111").Value = -4901
$ws.Range("N111").Value = -48098.333

# Hunk 9: sheet ALC, row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1378.25
$ws.Range("I132").Value = 1253.0385
$ws.Range("K132").Value = 3759.1155
$ws.Range("M132").Value = -1229.1155

# Hunk 10: sheet ALC, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3909.3225
$ws.Range("I138").Value = 2504.2
$ws.Range("J138").Value = 4578.4287
$ws.Range("K138").Value = 7512.599999999999
$ws.Range("L138").Value = 13735.2861
$ws.Range("M138").Value = -2372.599999999999
$ws.Range("N138").Value = -24015.2861

# Hunk 11: sheet ARM, row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1795.1333
$ws.Range("I2").Value = 1446.7273
$ws.Range("J2").Value = 2753.25
$ws.Range("K2").Value = 1446.7273
$ws.Range("L2").Value = 2753.25
$ws.Range("M2").Value = -1333.7273
$ws.Range("N2").Value = -2979.25

# Hunk 12: sheet ARM, row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1238
$ws.Range("I45").Value = 1272.5
$ws.Range("J45").Value = 1100
$ws.Range("K45").Value = 1272.5
$ws.Range("L45").Value = 1100
$ws.Range("M45").Value = -895.5
$ws.Range("N45").Value = -1854

# Hunk 13: sheet ARM, row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1795.1333
$ws.Range("I116").Value = 1446.7273
$ws.Range("J116").Value = 2753.25
$ws.Range("K116").Value = 1446.7273
$ws.Range("L116").Value = 2753.25
$ws.Range("M116").Value = 847.2727
$ws.Range("N116").Value = -7341.25

# Hunk 14: sheet ARM, row 123
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H123").Value = 35429
$ws.Range("J123").Value = 35429
$ws.Range("L123").Value = 35429
$ws.Range("N123").Value = -45229

# Hunk 15: sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 6341.2666
$ws.Range("I132").Value = 7905.8887
$ws.Range("J132").Value = 3994.3333
$ws.Range("K132").Value = 23717.6661
$ws.Range("L132").Value = 11982.9999
$ws.Range("M132").Value = -21187.6661
$ws.Range("N132").Value = -17042.9999

# Hunk 16: sheet BSM, row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1795.1333
$ws.Range("I3").Value = 1446.7273
$ws.Range("J3").Value = 2753.25
$ws.Range("K3").Value = 1446.7273
$ws.Range("L3").Value = 2753.25
$ws.Range("M3").Value = -1332.7273
$ws.Range("N3").Value = -2981.25

# Hunk 17: sheet BSM, row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1392.7646
$ws.Range("I94").Value = 1154.3636
$ws.Range("J94").Value = 1829.8334
$ws.Range("K94").Value = 1154.3636
$ws.Range("L94").Value = 1829.8334
$ws.Range("M94").Value = -703.3635999999999
$ws.Range("N94").Value = -2731.8334

# Hunk 18: sheet BSM, row 95
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 38986.75
$ws.Range("J95").Value = 38986.75
$ws.Range("L95").Value = 38986.75
$ws.Range("N95").Value = -44478.75

# Hunk 19: sheet BSM, row 100
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H100").Value = 22400
$ws.Range("J100").Value = 22400
$ws.Range("L100").Value = 22400
$ws.Range("N100").Value = -24564

# Hunk 20: sheet BSM, row 103
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 100000
$ws.Range("J103").Value = 100000
$ws.Range("L103").Value = 100000
$ws.Range("N103").Value = -102344

# Hunk 21: sheet BSM, row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3501.25
$ws.Range("I107").Value = 3987.5715
$ws.Range("J107").Value = 2820.4
$ws.Range("K107").Value = 3987.5715
$ws.Range("L107").Value = 2820.4
$ws.Range("M107").Value = -2067.5715
$ws.Range("N107").Value = -6660.4

# Hunk 22: sheet CRP, row 3
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 34999.75
$ws.Range("I3").Value = 49999.5
$ws.Range("J3").Value = 20000
$ws.Range("K3").Value = 49999.5
$ws.Range("L3").Value = 20000
$ws.Range("M3").Value = -49886.5
$ws.Range("N3").Value = -20226

# Hunk 23: sheet CRP, row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 447.6
$ws.Range("I22").Value = 243.38461
$ws.Range("K22").Value = 243.38461
$ws.Range("M22").Value = 106.61539

# Hunk 24: sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2489.261
$ws.Range("I31").Value = 1336.8485
$ws.Range("K31").Value = 1336.8485
$ws.Range("M31").Value = -1041.8485

# Hunk 25: sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2489.261
$ws.Range("I34").Value = 1336.8485
$ws.Range("K34").Value = 1336.8485
$ws.Range("M34").Value = -1134.8485

# Hunk 26: sheet CRP, row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2804.111
$ws.Range("I58").Value = 2563.6365
$ws.Range("J58").Value = 3182
$ws.Range("K58").Value = 2563.6365
$ws.Range("L58").Value = 3182
$ws.Range("M58").Value = -2360.6365
$ws.Range("N58").Value = -3588

# Hunk 27: sheet CRP, row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3214.7334
$ws.Range("I132").Value = 2482.2
$ws.Range("J132").Value = 4679.8
$ws.Range("K132").Value = 7446.599999999999
$ws.Range("L132").Value = 14039.4
$ws.Range("M132").Value = -4916.599999999999
$ws.Range("N132").Value = -19099.4

# Hunk 28: sheet CRP, row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2804.111
$ws.Range("I136").Value = 2563.6365
$ws.Range("J136").Value = 3182
$ws.Range("K136").Value = 7690.9095
$ws.Range("L136").Value = 9546
$ws.Range("M136").Value = -5140.9095
$ws.Range("N136").Value = -14646

# Hunk 29: sheet CUL, row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 7268
$ws.Range("I5").Value = 7268
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 21804
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -21692
$ws.Range("N5").ClearContents()

# Hunk 30: sheet CUL, row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1033.6666
$ws.Range("J122").Value = 1451
$ws.Range("L122").Value = 13059
$ws.Range("N122").Value = -17959

# Hunk 31: sheet CUL, row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 26317050
$ws.Range("I131").Value = 473
$ws.Range("J131").Value = 35715828
$ws.Range("K131").Value = 1419
$ws.Range("L131").Value = 107147484
$ws.Range("M131").Value = 3621
$ws.Range("N131").Value = -107157564

# Hunk 32: sheet CUL, row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2097.7334
$ws.Range("I132").Value = 1499
$ws.Range("K132").Value = 13491
$ws.Range("M132").Value = -10961

# Hunk 33: sheet CUL, row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 7268
$ws.Range("I135").Value = 7268
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 65412
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -62877
$ws.Range("N135").ClearContents()

# Hunk 34: sheet GSM, row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4039.3076
$ws.Range("I80").Value = 4000.7144
$ws.Range("K80").Value = 4000.7144
$ws.Range("M80").Value = -3002.7144

# Hunk 35: sheet GSM, row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 4039.3076
$ws.Range("I83").Value = 4000.7144
$ws.Range("K83").Value = 20003.572
$ws.Range("M83").Value = -15011.572

# Hunk 36: sheet LTW, row 17
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 10201.6
$ws.Range("I17").Value = 504
$ws.Range("K17").Value = 504
$ws.Range("M17").Value = -334

# Hunk 37: sheet LTW, row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2881
$ws.Range("I82").Value = 1948
$ws.Range("J82").Value = 4000.6
$ws.Range("K82").Value = 1948
$ws.Range("L82").Value = 4000.6
$ws.Range("M82").Value = -1587
$ws.Range("N82").Value = -4722.6

# Hunk 38: sheet LTW, row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2881
$ws.Range("I85").Value = 1948
$ws.Range("J85").Value = 4000.6
$ws.Range("K85").Value = 1948
$ws.Range("L85").Value = 4000.6
$ws.Range("M85").Value = -700
$ws.Range("N85").Value = -6496.6

# Hunk 39: sheet LTW, row 107
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H107").Value = 3500
$ws.Range("I107").Value = 3500
$ws.Range("K107").Value = 3500
$ws.Range("M107").Value = -1580

# Hunk 40: sheet LTW, row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2896.6538
$ws.Range("I136").Value = 2186.7693
$ws.Range("K136").Value = 6560.3079
$ws.Range("M136").Value = -4010.3079

# Hunk 41: sheet WVR, row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2076.7073
$ws.Range("I136").Value = 1883.7858
$ws.Range("J136").Value = 2492.2307
$ws.Range("K136").Value = 5651.357400000001
$ws.Range("L136").Value = 7476.6921
$ws.Range("M136").Value = -3101.357400000001
$ws.Range("N136").Value = -12576.6921
